$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$wsSchedule.Range("E2").Value = -31.38044325000001
$wsSchedule.Range("F2").Value = -0.6918087136243388
$wsSchedule.Range("E3").Value = 369.42665175
$wsSchedule.Range("F3").Value = 24.43297961309524

# --- Detailed sheet updates ---
$wsDetailed.Range("B13").Value = 57.06003
$wsDetailed.Range("B14").Value = 56.98
$wsDetailed.Range("C15").Value = "historical"
$wsDetailed.Range("C16").Value = "historical"
$wsDetailed.Range("B17").Value = 15.56016
$wsDetailed.Range("B18").Value = 0
$wsDetailed.Range("B19").Value = 0.7
$wsDetailed.Range("B20").Value = -4.52229
$wsDetailed.Range("B21").Value = -4.82627
$wsDetailed.Range("B22").Value = -5.50985
$wsDetailed.Range("B23").Value = 0.51
$wsDetailed.Range("B24").Value = -0.95265
$wsDetailed.Range("B25").Value = 0
$wsDetailed.Range("B26").Value = -0.90081
$wsDetailed.Range("B27").Value = -0.88996
$wsDetailed.Range("B28").Value = -5.37499
$wsDetailed.Range("B29").Value = -6.12399
$wsDetailed.Range("B30").Value = -6.74277
$wsDetailed.Range("B31").Value = -20.59074
$wsDetailed.Range("B32").Value = -14.45981
$wsDetailed.Range("B33").Value = -11.01
$wsDetailed.Range("B34").Value = -7.1669
$wsDetailed.Range("B35").Value = -6.57264
$wsDetailed.Range("B37").Value = 0.66204
$wsDetailed.Range("B38").Value = 3.98304
$wsDetailed.Range("B39").Value = 11.98336
$wsDetailed.Range("B40").Value = 40.9994
$wsDetailed.Range("B41").Value = 56.98
$wsDetailed.Range("B43").Value = 56.98
$wsDetailed.Range("B44").Value = 56.98
$wsDetailed.Range("B45").Value = 43.56411
$wsDetailed.Range("B46").Value = 55.89482
